# ---------------------------------------------------------------------------
# "updated word count on title page after test run"
#
# 1. Remove the stray empty run (<w:r><w:t xml:space="preserve"/></w:r>) left
#    behind in the trailing FirstParagraph-styled paragraph at the end of the
#    document body - it has no visible text, it's just noise from a prior
#    save/word-count pass.
# 2. Re-base the built-in "Subtitle" style on "Title" instead of "Normal",
#    and strip the explicit font color override from its rPr (the Subtitle
#    Char linked character style keeps its own color untouched).
# 3. Strip the explicit font color override from the "AbstractTitle" rPr.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. Drop the empty trailing run -----------------------------------------
$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range
# Trim the paragraph mark off the end so we only touch the run(s) that
# precede it; this removes the empty run while leaving the (empty) paragraph
# itself - and its pStyle - intact. (The run carries no visible text, so the
# trimmed range is collapsed/zero-length by character count, but Delete()
# still removes the underlying empty <w:r> element.)
$trimmed = $d.Range($lastRange.Start, $lastRange.End - 1)
$trimmed.Delete()

# --- 2. Subtitle style: re-parent to Title, drop its own color override ----
$subtitle = $d.Styles("Subtitle")
$subtitle.BaseStyle = $d.Styles("Title")
$subtitle.Font.Color = -16777216   ; # wdColorAutomatic - clears the explicit override

# --- 3. AbstractTitle style: drop its color override ------------------------
$abstractTitle = $d.Styles("AbstractTitle")
$abstractTitle.Font.Color = -16777216   ; # wdColorAutomatic
